$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.706.68"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "3.166.17"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'207.15"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'609.64"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +4.98%  "
$ws.Range("D8").Value = "'0.672"
$ws.Range("E8").Value = "  +4.89%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.164.58"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("D11").Value = "'0.534"
$ws.Range("E11").Value = "  -6.62%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").Value = "'5.27"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "3.746.24"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "87.522.47"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "'32.27"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "3.171.08"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").Value = "'3.18"
$ws.Range("E19").Value = "  +8.14%  "
$ws.Range("D20").Value = "'13.38"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").Value = "'410.37"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").Value = "'8.42"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").Value = "'5.06"
$ws.Range("E23").Value = "  -4.51%  "
$ws.Range("D24").Value = "'5.21"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").Value = "'12.22"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "3.334.87"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("E27").Value = "  +5.44%  "
$ws.Range("D28").Value = "'73.31"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'0.165"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'547.41"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "'8.24"
$ws.Range("E33").Value = "  -5.81%  "
$ws.Range("D34").Value = "'1.31"
$ws.Range("E34").Value = "  -6.73%  "
$ws.Range("D35").Value = "'6.82"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'1.85"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("D37").Value = "'0.130"
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("D38").Value = "'21.80"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'21.80"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "'3.03"
$ws.Range("E41").Value = "  +5.31%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'1.90"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").Value = "'0.371"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "'151.09"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "'173.17"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'0.124"
$ws.Range("E48").Value = "  +5.51%  "
$ws.Range("E49").Value = "  -6.68%  "
$ws.Range("D50").Value = "'23.96"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  -6.08%  "
